$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Fecha (D2) and Volumen (M2)
$ws.Range("D2").Value = 44250
$ws.Range("M2").Value = 200

# Row 3: Fecha (D3) and Volumen (M3)
$ws.Range("D3").Value = 44257
$ws.Range("M3").Value = 100

# Row 4: Fecha (D4) and Volumen (M4)
$ws.Range("D4").Value = 44253
$ws.Range("M4").Value = 160
